$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1943.9445
$ws.Range("J17").Value = 1943.9445
$ws.Range("L17").Value = 5831.833500000001
$ws.Range("N17").Value = -6167.833500000001
$ws.Range("H41").Value = 949
$ws.Range("I41").Value = 1000.5
$ws.Range("J41").Value = 897.5
$ws.Range("K41").Value = 1000.5
$ws.Range("L41").Value = 897.5
$ws.Range("M41").Value = -560.5
$ws.Range("N41").Value = -1777.5
$ws.Range("H58").Value = 566.9
$ws.Range("I58").Value = 204.375
$ws.Range("J58").Value = 2017
$ws.Range("K58").Value = 613.125
$ws.Range("L58").Value = 6051
$ws.Range("M58").Value = -463.125
$ws.Range("N58").Value = -6351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 29999.5
$ws.Range("I76").Value = 29999
$ws.Range("K76").Value = 29999
$ws.Range("M76").Value = -29661
$ws.Range("H79").Value = 29999.5
$ws.Range("I79").Value = 29999
$ws.Range("K79").Value = 29999
$ws.Range("M79").Value = -28829
$ws.Range("H122").Value = 1976.6471
$ws.Range("I122").Value = 1976.6471
$ws.Range("K122").Value = 5929.9413
$ws.Range("M122").Value = -3479.9413
$ws.Range("H131").Value = 73749.5
$ws.Range("J131").Value = 73749.5
$ws.Range("L131").Value = 73749.5
$ws.Range("N131").Value = -83829.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 20798.125
$ws.Range("J88").Value = 20798.125
$ws.Range("L88").Value = 20798.125
$ws.Range("N88").Value = -21610.125
$ws.Range("H91").Value = 20798.125
$ws.Range("J91").Value = 20798.125
$ws.Range("L91").Value = 20798.125
$ws.Range("N91").Value = -23606.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 207.11111
$ws.Range("I7").Value = 193.4
$ws.Range("J7").Value = 246.28572
$ws.Range("K7").Value = 193.4
$ws.Range("L7").Value = 246.28572
$ws.Range("M7").Value = -80.40000000000001
$ws.Range("N7").Value = -472.28572
$ws.Range("H22").Value = 3333782
$ws.Range("I22").Value = 577.5
$ws.Range("J22").Value = 5000384.5
$ws.Range("K22").Value = 577.5
$ws.Range("L22").Value = 5000384.5
$ws.Range("M22").Value = -227.5
$ws.Range("N22").Value = -5001084.5
$ws.Range("H99").Value = 2631.5
$ws.Range("I99").Value = 2631.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2631.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1133.5
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 2631.5
$ws.Range("I126").Value = 2631.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7894.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5424.5
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1343.625
$ws.Range("I5").Value = 1422.5
$ws.Range("J5").Value = 1264.75
$ws.Range("K5").Value = 4267.5
$ws.Range("L5").Value = 3794.25
$ws.Range("M5").Value = -4155.5
$ws.Range("N5").Value = -4018.25
$ws.Range("H34").Value = 62409.777
$ws.Range("J34").Value = 64316.234
$ws.Range("L34").Value = 192948.702
$ws.Range("N34").Value = -193116.702
$ws.Range("H39").Value = 6278.4443
$ws.Range("J39").Value = 6687.875
$ws.Range("L39").Value = 20063.625
$ws.Range("N39").Value = -20651.625
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 800
$ws.Range("K49").Value = 2400
$ws.Range("M49").Value = -2244
$ws.Range("H55").Value = 8299
$ws.Range("J55").Value = 9798
$ws.Range("L55").Value = 29394
$ws.Range("N55").Value = -29748
$ws.Range("H80").Value = 5806.5
$ws.Range("I80").Value = 5708.8335
$ws.Range("K80").Value = 17126.5005
$ws.Range("M80").Value = -16190.5005
$ws.Range("H83").Value = 5806.5
$ws.Range("I83").Value = 5708.8335
$ws.Range("K83").Value = 51379.5015
$ws.Range("M83").Value = -46699.5015
$ws.Range("H87").Value = 4997.5
$ws.Range("I87").Value = 4997.5
$ws.Range("K87").Value = 14992.5
$ws.Range("M87").Value = -13744.5
$ws.Range("H88").Value = 21000
$ws.Range("J88").Value = 21000
$ws.Range("L88").Value = 63000
$ws.Range("N88").Value = -63856
$ws.Range("H90").Value = 4997.5
$ws.Range("I90").Value = 4997.5
$ws.Range("K90").Value = 44977.5
$ws.Range("M90").Value = -38737.5
$ws.Range("H91").Value = 21000
$ws.Range("J91").Value = 21000
$ws.Range("L91").Value = 63000
$ws.Range("N91").Value = -65964
$ws.Range("H92").Value = 286.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 286.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 859.5
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3355.5
$ws.Range("H107").Value = 734
$ws.Range("I107").Value = 751
$ws.Range("J107").Value = 674.5
$ws.Range("K107").Value = 2253
$ws.Range("L107").Value = 2023.5
$ws.Range("M107").Value = -333
$ws.Range("N107").Value = -5863.5
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 1694.5385
$ws.Range("J132").Value = 1833
$ws.Range("L132").Value = 16497
$ws.Range("N132").Value = -21557
$ws.Range("H135").Value = 1343.625
$ws.Range("I135").Value = 1422.5
$ws.Range("J135").Value = 1264.75
$ws.Range("K135").Value = 12802.5
$ws.Range("L135").Value = 11382.75
$ws.Range("M135").Value = -10267.5
$ws.Range("N135").Value = -16452.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 91799.39999999999
$ws.Range("I70").Value = 89749.5
$ws.Range("K70").Value = 89749.5
$ws.Range("M70").Value = -89479.5
$ws.Range("H73").Value = 91799.39999999999
$ws.Range("I73").Value = 89749.5
$ws.Range("K73").Value = 89749.5
$ws.Range("M73").Value = -88813.5
$ws.Range("H80").Value = 2408.4
$ws.Range("I80").Value = 2399
$ws.Range("K80").Value = 2399
$ws.Range("M80").Value = -1401
$ws.Range("H83").Value = 2408.4
$ws.Range("I83").Value = 2399
$ws.Range("K83").Value = 11995
$ws.Range("M83").Value = -7003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 80715
$ws.Range("J36").Value = 80715
$ws.Range("L36").Value = 80715
$ws.Range("N36").Value = -81839
$ws.Range("H130").Value = 66657.5
$ws.Range("J130").Value = 66657.5
$ws.Range("L130").Value = 66657.5
$ws.Range("N130").Value = -76697.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11620.75
$ws.Range("I62").Value = 4998.5
$ws.Range("J62").Value = 13828.167
$ws.Range("K62").Value = 4998.5
$ws.Range("L62").Value = 13828.167
$ws.Range("M62").Value = -4374.5
$ws.Range("N62").Value = -15076.167
$ws.Range("H65").Value = 11620.75
$ws.Range("I65").Value = 4998.5
$ws.Range("J65").Value = 13828.167
$ws.Range("K65").Value = 24992.5
$ws.Range("L65").Value = 69140.83499999999
$ws.Range("M65").Value = -21872.5
$ws.Range("N65").Value = -75380.83499999999
$ws.Range("H122").Value = 3283.913
$ws.Range("I122").Value = 2823.2104
$ws.Range("K122").Value = 8469.6312
$ws.Range("M122").Value = -6019.6312
$ws.Range("H130").Value = 37332.668
$ws.Range("J130").Value = 37332.668
$ws.Range("L130").Value = 37332.668
$ws.Range("N130").Value = -47372.668

Write-Host "Applied 192 cell updates across 8 sheets"